# Remove three obsolete security-rule rows and the stray hyperlink that
# pointed to http://www.google.com/ (CodeQuality-rules-latest.xlsx cleanup).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the hyperlink (and its visible text) that lives on the row for
# "CQRules:CQBP-44---ConsecutivelyLogAndThrow" (column F), before the row
# indices shift due to the row deletions below.
[void]$ws.Range("F103").Hyperlinks.Delete()
[void]$ws.Range("F103").ClearContents()

# Delete the rows for squid:S3318, squid:S2078 and squid:S2076 (in
# descending row-number order so earlier row numbers stay valid while
# deleting).
[void]$ws.Rows.Item(26).Delete()
[void]$ws.Rows.Item(10).Delete()
[void]$ws.Rows.Item(8).Delete()

# Leave the selection where the author ended up after clearing the
# hyperlink cell.
[void]$ws.Range("F100").Select()
